$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecasting")

# The FX forecast curve creation now takes a single "currencyPair" (e.g. "USDZAR")
# instead of separate "baseCurrency"/"counterCurrency" values (SerializableViaName
# converters collapse the pair into one name). This removes one input row and
# shifts the remaining parameter rows up by one.

$ws.Range("B16").Value = "currencyPair"
$ws.Range("C16").Value = "USDZAR"

$ws.Range("B17").Value = "fxRateAtAnchorDate"
$ws.Range("C17").Value = 13.66

$ws.Range("B18").Value = "baseCurrencyFXBasisCurve"
$ws.Range("C18").Formula = "=E12"

$ws.Range("B19").Value = "counterCurrencyFXBasisCurve"
$ws.Range("C19").Formula = "=B12"

$ws.Range("B20:C20").Clear()

$ws.Range("B22").Formula = "=_xll.QSA.CreateFXForecastCurve(C15,C16,C17,C18,C19)"
